$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report week date range) ---
$ws.Range("A8").Value = "Volume 32   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/1/2025  Through  12/7/2025"

# --- Cells converting from placeholder text to numeric values ---
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = '#,##0'
$ws.Range("H15").Value = -100
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D28").Value = 2
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = 50
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C33").Value = 1
$ws.Range("C33").NumberFormat = '#,##0'
$ws.Range("F33").Value = 1
$ws.Range("F33").NumberFormat = '#,##0'

# --- Cells converting from numeric values to placeholder text "0" (style matches D14) ---
$ws.Range("D14").Copy($ws.Range("C17"))
$ws.Range("D14").Copy($ws.Range("C22"))
$ws.Range("D14").Copy($ws.Range("C23"))

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("J15").Value = 16
$ws.Range("K15").Value = -31.25
$ws.Range("N15").Value = -63.333333333333
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 158
$ws.Range("J16").Value = 196
$ws.Range("K16").Value = -19.387755102040
$ws.Range("L16").Value = -21.393034825870
$ws.Range("M16").Value = 18.796992481203
$ws.Range("N16").Value = -87.770897832817
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 7.692307692307
$ws.Range("I17").Value = 198
$ws.Range("J17").Value = 185
$ws.Range("K17").Value = 7.027027027027
$ws.Range("L17").Value = 13.793103448275
$ws.Range("M17").Value = 122.47191011236
$ws.Range("N17").Value = -28.776978417266
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -37.5
$ws.Range("I18").Value = 261
$ws.Range("J18").Value = 230
$ws.Range("K18").Value = 13.478260869565
$ws.Range("L18").Value = 6.097560975609
$ws.Range("M18").Value = 14.977973568281
$ws.Range("N18").Value = -90.757790368272
$ws.Range("C19").Value = 31
$ws.Range("D19").Value = 39
$ws.Range("E19").Value = -20.512820512820
$ws.Range("F19").Value = 107
$ws.Range("G19").Value = 144
$ws.Range("H19").Value = -25.694444444444
$ws.Range("I19").Value = 1513
$ws.Range("J19").Value = 1561
$ws.Range("K19").Value = -3.074951953875
$ws.Range("L19").Value = -5.732087227414
$ws.Range("M19").Value = 28.656462585034
$ws.Range("N19").Value = -57.078014184397
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 400
$ws.Range("I20").Value = 58
$ws.Range("K20").Value = -27.5
$ws.Range("L20").Value = -60.544217687074
$ws.Range("M20").Value = -31.764705882352
$ws.Range("N20").Value = -98.216482164821
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 50
$ws.Range("E21").Value = -24
$ws.Range("F21").Value = 156
$ws.Range("G21").Value = 199
$ws.Range("H21").Value = -21.608040201005
$ws.Range("I21").Value = 2199
$ws.Range("J21").Value = 2270
$ws.Range("K21").Value = -3.127753303964
$ws.Range("L21").Value = -7.914572864321
$ws.Range("M21").Value = 27.478260869565
$ws.Range("N21").Value = -80.376583972871
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("L22").Value = -35.714285714285
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 50
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 84
$ws.Range("E24").Value = -52.380952380952
$ws.Range("F24").Value = 166
$ws.Range("G24").Value = 295
$ws.Range("H24").Value = -43.728813559322
$ws.Range("I24").Value = 2633
$ws.Range("J24").Value = 3166
$ws.Range("K24").Value = -16.835123183828
$ws.Range("L24").Value = -9.456671251719
$ws.Range("M24").Value = 66.751108296390
$ws.Range("C25").Value = 28
$ws.Range("D25").Value = 77
$ws.Range("E25").Value = -63.636363636363
$ws.Range("F25").Value = 130
$ws.Range("G25").Value = 272
$ws.Range("H25").Value = -52.205882352941
$ws.Range("I25").Value = 2162
$ws.Range("J25").Value = 2783
$ws.Range("K25").Value = -22.314049586776
$ws.Range("L25").Value = -14.747634069400
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -62.5
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = -34.285714285714
$ws.Range("I26").Value = 336
$ws.Range("J26").Value = 350
$ws.Range("K26").Value = -4
$ws.Range("L26").Value = 1.510574018126
$ws.Range("M26").Value = 3.067484662576
$ws.Range("G27").Value = 2
$ws.Range("J27").Value = 22
$ws.Range("K27").Value = -36.363636363636
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 28.571428571428
$ws.Range("I28").Value = 94
$ws.Range("J28").Value = 104
$ws.Range("K28").Value = -9.615384615384
$ws.Range("L28").Value = 13.253012048192
$ws.Range("L31").Value = -42.307692307692
$ws.Range("I33").Value = 4
$ws.Range("K33").Value = 33.333333333333
$ws.Range("L33").Value = 33.333333333333
